# "Added last minute updates"
#
# Target paragraph is the very first paragraph of the body, which holds
# the hidden opendope bookmark id used for topic selection. The edit:
#   1. adds a 4-sided paragraph border with 5-twip spacing (pBdr),
#   2. bumps the left indent from 120 -> 225 twips,
#   3. renames the bookmark id text and drops the trailing run that only
#      held a single space.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# 1) Paragraph border, 5 twips of space on every side -> <w:pBdr><w:top
#    w:space="5"/><w:left .../><w:bottom .../><w:right .../></w:pBdr>
$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# 2) Left indent 120 twips (6pt) -> 225 twips (11.25pt)
$p.Format.LeftIndent = 11.25

# 3) Replace the whole paragraph's text (both runs, but not the paragraph
#    mark) with the new bookmark id. This merges the two existing runs
#    (the id run + the trailing " " run) into a single run carrying the
#    first run's formatting, matching the diff which drops the
#    space-only run entirely.
$full = $p.Range
$body = $d.Range($full.Start, $full.End - 1)
$body.Text = "**ID__AFFARS_AFMC_PGI_5315_3__ID**"
